# Applies the Oliva (1992) stress-clash generator/evaluator update.
#
# The generated example list (col A, "Exemple") is re-derived with a new
# ordering, and the per-example metrics are recomputed: RC2 (col C), RC3
# (col D), Complexitat = RC2+RC3 (col E), and the 1-indexed odd-syllable
# positions where each clash type occurs (col F = "RC2 posicions", col G =
# "RC3 posicions", comma-separated text such as "1, 3").
#
# Only cells whose value actually changes are touched (155 of the 7x33
# grid), so already-blank "posicions" cells that stay blank are left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "RC2 posicions" / "RC3 posicions" (F, G) hold digit lists ("1, 3") as
# text - force text formatting on those columns first so Excel doesn't
# auto-coerce a lone-digit value (e.g. "5") into a numeric cell.
$ws.Range("F2:G33").NumberFormat = "@"

$changes = @(
    ,@(2, 1, 'ATATATATAT', 't')
    ,@(2, 3, 0, 'n')
    ,@(2, 5, 0, 'n')
    ,@(2, 6, '', 't')
    ,@(3, 1, 'ATAAATATAT', 't')
    ,@(3, 3, 1, 'n')
    ,@(3, 4, 0, 'n')
    ,@(3, 6, '3', 't')
    ,@(3, 7, '', 't')
    ,@(4, 1, 'ATATTAATAT', 't')
    ,@(4, 7, '5', 't')
    ,@(5, 1, 'ATATAAATAT', 't')
    ,@(5, 3, 1, 'n')
    ,@(5, 4, 0, 'n')
    ,@(5, 5, 1, 'n')
    ,@(5, 6, '5', 't')
    ,@(5, 7, '', 't')
    ,@(6, 1, 'TAATATATAT', 't')
    ,@(6, 3, 0, 'n')
    ,@(6, 5, 1, 'n')
    ,@(6, 6, '', 't')
    ,@(6, 7, '1', 't')
    ,@(7, 1, 'ATTAATATAT', 't')
    ,@(7, 3, 0, 'n')
    ,@(7, 4, 1, 'n')
    ,@(7, 5, 1, 'n')
    ,@(7, 6, '', 't')
    ,@(7, 7, '3', 't')
    ,@(8, 1, 'AAATATATAT', 't')
    ,@(8, 4, 0, 'n')
    ,@(8, 5, 1, 'n')
    ,@(8, 6, '1', 't')
    ,@(8, 7, '', 't')
    ,@(9, 1, 'ATATATTAAT', 't')
    ,@(9, 7, '7', 't')
    ,@(10, 1, 'ATATATAAAT', 't')
    ,@(10, 4, 0, 'n')
    ,@(10, 5, 1, 'n')
    ,@(10, 6, '7', 't')
    ,@(10, 7, '', 't')
    ,@(11, 1, 'ATATAAAAAT', 't')
    ,@(11, 3, 2, 'n')
    ,@(11, 5, 2, 'n')
    ,@(11, 6, '5, 7', 't')
    ,@(12, 1, 'ATATAATAAT', 't')
    ,@(12, 3, 1, 'n')
    ,@(12, 4, 1, 'n')
    ,@(12, 6, '5', 't')
    ,@(12, 7, '7', 't')
    ,@(13, 1, 'AATAATATAT', 't')
    ,@(13, 3, 1, 'n')
    ,@(13, 4, 1, 'n')
    ,@(13, 6, '1', 't')
    ,@(13, 7, '3', 't')
    ,@(14, 1, 'TAATATAAAT', 't')
    ,@(14, 3, 1, 'n')
    ,@(14, 4, 1, 'n')
    ,@(14, 5, 2, 'n')
    ,@(14, 6, '7', 't')
    ,@(14, 7, '1', 't')
    ,@(15, 1, 'TAAAATATAT', 't')
    ,@(15, 3, 1, 'n')
    ,@(15, 4, 1, 'n')
    ,@(15, 6, '3', 't')
    ,@(15, 7, '1', 't')
    ,@(16, 1, 'ATAAAAATAT', 't')
    ,@(16, 6, '3, 5', 't')
    ,@(17, 1, 'ATAATAATAT', 't')
    ,@(17, 3, 1, 'n')
    ,@(17, 4, 1, 'n')
    ,@(17, 5, 2, 'n')
    ,@(17, 6, '3', 't')
    ,@(17, 7, '5', 't')
    ,@(18, 1, 'ATTAATAAAT', 't')
    ,@(18, 3, 1, 'n')
    ,@(18, 5, 2, 'n')
    ,@(18, 6, '7', 't')
    ,@(18, 7, '3', 't')
    ,@(19, 1, 'AAATTAATAT', 't')
    ,@(19, 3, 1, 'n')
    ,@(19, 4, 1, 'n')
    ,@(19, 5, 2, 'n')
    ,@(19, 6, '1', 't')
    ,@(19, 7, '5', 't')
    ,@(20, 1, 'AAATATAAAT', 't')
    ,@(20, 3, 2, 'n')
    ,@(20, 5, 2, 'n')
    ,@(20, 6, '1, 7', 't')
    ,@(21, 1, 'AAAAATATAT', 't')
    ,@(21, 3, 2, 'n')
    ,@(21, 4, 0, 'n')
    ,@(21, 5, 2, 'n')
    ,@(21, 6, '1, 3', 't')
    ,@(21, 7, '', 't')
    ,@(22, 1, 'AAATAAATAT', 't')
    ,@(22, 3, 2, 'n')
    ,@(22, 4, 0, 'n')
    ,@(22, 6, '1, 5', 't')
    ,@(22, 7, '', 't')
    ,@(23, 1, 'ATAAATAAAT', 't')
    ,@(23, 3, 2, 'n')
    ,@(23, 5, 2, 'n')
    ,@(23, 6, '3, 7', 't')
    ,@(24, 1, 'AAATATTAAT', 't')
    ,@(24, 3, 1, 'n')
    ,@(24, 4, 1, 'n')
    ,@(24, 5, 2, 'n')
    ,@(24, 6, '1', 't')
    ,@(24, 7, '7', 't')
    ,@(25, 1, 'ATAAAAAAAT', 't')
    ,@(25, 3, 3, 'n')
    ,@(25, 4, 0, 'n')
    ,@(25, 5, 3, 'n')
    ,@(25, 6, '3, 5, 7', 't')
    ,@(25, 7, '', 't')
    ,@(26, 1, 'AAAAATAAAT', 't')
    ,@(26, 3, 3, 'n')
    ,@(26, 5, 3, 'n')
    ,@(26, 6, '1, 3, 7', 't')
    ,@(27, 1, 'AAATAATAAT', 't')
    ,@(27, 3, 2, 'n')
    ,@(27, 5, 3, 'n')
    ,@(27, 6, '1, 5', 't')
    ,@(27, 7, '7', 't')
    ,@(28, 1, 'AATAATAAAT', 't')
    ,@(28, 3, 2, 'n')
    ,@(28, 4, 1, 'n')
    ,@(28, 5, 3, 'n')
    ,@(28, 6, '1, 7', 't')
    ,@(28, 7, '3', 't')
    ,@(29, 1, 'AAAATAATAT', 't')
    ,@(29, 6, '1, 3', 't')
    ,@(29, 7, '5', 't')
    ,@(30, 1, 'AAATAAAAAT', 't')
    ,@(30, 3, 3, 'n')
    ,@(30, 4, 0, 'n')
    ,@(30, 6, '1, 5, 7', 't')
    ,@(30, 7, '', 't')
    ,@(31, 1, 'TAAAATAAAT', 't')
    ,@(31, 4, 1, 'n')
    ,@(31, 5, 3, 'n')
    ,@(31, 6, '3, 7', 't')
    ,@(31, 7, '1', 't')
    ,@(32, 1, 'AAAAAAATAT', 't')
    ,@(32, 3, 3, 'n')
    ,@(32, 4, 0, 'n')
    ,@(32, 5, 3, 'n')
    ,@(32, 6, '1, 3, 5', 't')
    ,@(32, 7, '', 't')
    ,@(33, 1, 'AAAAAAAAAT', 't')
    ,@(33, 3, 4, 'n')
    ,@(33, 4, 0, 'n')
    ,@(33, 5, 4, 'n')
    ,@(33, 6, '1, 3, 5, 7', 't')
    ,@(33, 7, '', 't')
)

foreach ($chg in $changes) {
    $r    = $chg[0]
    $c    = $chg[1]
    $val  = $chg[2]
    $kind = $chg[3]

    $cell = $ws.Cells.Item($r, $c)
    if ($kind -eq 'n') {
        $cell.Value = $val
    } else {
        $cell.Value = [string]$val
    }
}
